# Fix ordre virement (carte_sejour,cin) condition issue
# Replace the single data row with the corrected/expanded tax-state rows (A2:O9)
# and regenerate the totals row (now row 10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "908/LF/DIRECTION REGIONALE SUD"
$ws.Range("B2").Value = "Logement de fonction"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "354646"
$ws.Range("D2").Value = "AGENCE KHATABI"
$ws.Range("E2").Value = "oui"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 8000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 7200

# Row 3
$ws.Range("A3").Value = "908/LF/DIRECTION REGIONALE SUD"
$ws.Range("B3").Value = "Logement de fonction"
$ws.Range("C3").Value = "B12346"
$ws.Range("D3").Value = "BAKKALI MOHAMED"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 6000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "--"
$ws.Range("O3").Value = 5400

# Row 4
$ws.Range("A4").Value = "908/LF/DIRECTION REGIONALE SUD"
$ws.Range("B4").Value = "Logement de fonction"
$ws.Range("C4").Value = "L234567"
$ws.Range("D4").Value = "NACER YASSINE"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "--"
$ws.Range("O4").Value = 2000

# Row 5
$ws.Range("A5").Value = "389/AOURIR"
$ws.Range("B5").Value = "Point de vente"
$ws.Range("C5").Value = "BJ36877"
$ws.Range("D5").Value = "CHARIJI ABDELLAH"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 7000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = "--"
$ws.Range("O5").Value = 6300

# Row 6
$ws.Range("A6").Value = "908/DIRECTION REGIONALE SUD"
$ws.Range("B6").Value = "Direction régionale"
$ws.Range("C6").Value = "J207703"
$ws.Range("D6").Value = "ACHENGLI LAILA"
$ws.Range("E6").Value = "non"
$ws.Range("F6").Value = "mensuelle"
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = 20000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = "--"
$ws.Range("O6").Value = 17000

# Row 7
$ws.Range("A7").Value = "900/PATIO"
$ws.Range("B7").Value = "Siège"
$ws.Range("C7").Value = "J207703"
$ws.Range("D7").Value = "ACHENGLI LAILA"
$ws.Range("E7").Value = "non"
$ws.Range("F7").Value = "mensuelle"
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 4500
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = "--"
$ws.Range("O7").Value = 4050

# Row 8
$ws.Range("A8").Value = "001/SUP SUD"
$ws.Range("B8").Value = "Supervision"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "354646"
$ws.Range("D8").Value = "AGENCE KHATABI"
$ws.Range("E8").Value = "oui"
$ws.Range("F8").Value = "mensuelle"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 2400
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = "--"
$ws.Range("O8").Value = 2400

# Row 9
$ws.Range("A9").Value = "805/KOUTOUBIA"
$ws.Range("B9").Value = "Point de vente"
$ws.Range("C9").Value = "L234567"
$ws.Range("D9").Value = "NACER YASSINE"
$ws.Range("E9").Value = "non"
$ws.Range("F9").Value = "mensuelle"
$ws.Range("G9").Value = 15
$ws.Range("H9").Value = 12000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1800
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = "--"
$ws.Range("O9").Value = 10200

# Row 10
$ws.Range("A10").Value = " "
$ws.Range("B10").Value = " "
$ws.Range("C10").Value = " "
$ws.Range("D10").Value = " "
$ws.Range("E10").Value = " "
$ws.Range("F10").Value = " "
$ws.Range("G10").Value = " "
$ws.Range("H10").Value = 61900
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 7350
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 54550
